$wb = $excel.ActiveWorkbook

# --- 1. Rename Sheet1 -> "Binary Search 1" --------------------------------
$ws1 = $wb.Worksheets.Item(1)
$ws1.Name = "Binary Search 1"

# --- 2. Update the "Page No. in notes" column labels on rows 3-6 ----------
$ws1.Range("C3").Value = "Binary S1 10"
$ws1.Range("C4").Value = "Binary S1 12"
$ws1.Range("C5").Value = "Binary S1 13"
$ws1.Range("C6").Value = "Binary S1 14"

# --- 3. Insert new row 8 (Matrix median) -----------------------------------
$ws1.Range("B8").Value = 6
$ws1.Range("C8").Value = "Binary S1 18"
$ws1.Range("D8").Value = "Matrix median"

$ws1.Hyperlinks.Add($ws1.Range("E8"), "https://www.scaler.com/academy/mentee-dashboard/class/30364/homework/problems/357?navref=cl_tt_lst_sl", "", "", "https://www.scaler.com/academy/mentee-dashboard/class/30364/homework/problems/357?navref=cl_tt_lst_sl") | Out-Null
$ws1.Range("E5").Copy() | Out-Null
$ws1.Range("E8").PasteSpecial(-4122) | Out-Null

$ws1.Hyperlinks.Add($ws1.Range("F8"), "https://github.com/ankurnecessary/dsa/blob/main/1_binarySearch/6_matrix_median.java", "", "", "dsa/6_matrix_median.java at main · ankurnecessary/dsa · GitHub") | Out-Null
$ws1.Range("F4").Copy() | Out-Null
$ws1.Range("F8").PasteSpecial(-4122) | Out-Null

$ws1.Rows.Item(8).RowHeight = 72

# --- 4. Fix up C7's label (was "Binary S" -> "Binary S1 16") --------------
$ws1.Range("C7").Value = "Binary S1 16"

# --- 5. Update dimension / selection on sheet 1 ----------------------------
$ws1.Range("C3").Select() | Out-Null

# --- 6. Add "Binary Search 2" sheet right after "Binary Search 1" ---------
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "Binary Search 2"

$ws2.Range("B2").Value = "S.no."
$ws2.Range("C2").Value = "Page No. in notes"
$ws2.Range("D2").Value = "Question"
$ws2.Range("E2").Value = "Link"
$ws2.Range("F2").Value = "Github Link"

$ws2.Range("B3").Value = 1
$ws2.Range("C3").Value = "Binary S2 2"
$ws2.Range("D3").Value = "Square root of Integer"

$ws2.Hyperlinks.Add($ws2.Range("E3"), "https://www.scaler.com/academy/mentee-dashboard/class/30365/assignment/problems/200/?navref=cl_pb_nv_tb", "", "", "https://www.scaler.com/academy/mentee-dashboard/class/30365/assignment/problems/200/?navref=cl_pb_nv_tb") | Out-Null

$ws2.Range("B2:F3").Font.Name = "Calibri"
$ws2.Range("B2:F3").HorizontalAlignment = -4131
$ws2.Range("B2:F3").VerticalAlignment = -4160
$ws2.Range("E3").WrapText = $true

$ws2.Rows.Item(3).RowHeight = 72

$ws2.Range("F3").Select() | Out-Null

Write-Output "done"
